# 13. Monthly Budget.xlsx - add "MINIFS n MAXIFS" sheet
$wb = $excel.ActiveWorkbook

# --- Update the existing "Monthly Budget  (2)" totals row to use absolute
#     row references (B$3:B$7 etc.) instead of relative ones.
$ws2 = $wb.Worksheets.Item("Monthly Budget  (2)")
$ws2.Range("B8").Formula = "=SUM(B`$3:B`$7)"
$ws2.Range("C8").Formula = "=SUM(C`$3:C`$7)"
$ws2.Range("D8").Formula = "=SUM(D`$3:D`$7)"

# --- Add the new worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "MINIFS n MAXIFS"

# --- Header row ---
$ws3.Range("A1").Value = "Person"
$ws3.Range("B1").Value = "Year"
$ws3.Range("C1").Value = "Donation"

# --- Ayaan rows (2-5) ---
$ws3.Range("A2").Value = 500
$ws3.Range("B2").Value = "Ayaan"
$ws3.Range("A3").Value = 900
$ws3.Range("B3").Value = "Ayaan"
$ws3.Range("A4").Value = 5900
$ws3.Range("B4").Value = "Ayaan"
$ws3.Range("A5").Value = 700
$ws3.Range("B5").Value = "Ayaan"

$ws3.Range("D2").Value = "Ayaan Min"
$ws3.Range("E2").Formula = "=MINIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Ayaan"")"
$ws3.Range("D3").Value = "Ayaan Max"
$ws3.Range("E3").Formula = "=MAXIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Ayaan"")"
$ws3.Range("D4").Value = "Ayaan Sum"
$ws3.Range("E4").Formula = "=SUMIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Ayaan"")"

# --- Ayaan subtotal row (6) ---
$ws3.Range("A6").Formula = "=SUBTOTAL(9,A2:A5)"
$ws3.Range("B6").Value = "Ayaan Total"

# --- Naman rows (7-11) ---
$ws3.Range("A7").Value = 450
$ws3.Range("B7").Value = "Naman"
$ws3.Range("A8").Value = 200
$ws3.Range("B8").Value = "Naman"
$ws3.Range("A9").Value = 500
$ws3.Range("B9").Value = "Naman"
$ws3.Range("A10").Value = 90
$ws3.Range("B10").Value = "Naman"
$ws3.Range("A11").Value = 100
$ws3.Range("B11").Value = "Naman"

$ws3.Range("D7").Value = "Ayaan Min"
$ws3.Range("E7").Formula = "=MINIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Naman"")"
$ws3.Range("D8").Value = "Ayaan Max"
$ws3.Range("E8").Formula = "=MAXIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Naman"")"
$ws3.Range("D9").Value = "Ayaan Sum"
$ws3.Range("E9").Formula = "=SUMIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Naman"")"

# --- Naman subtotal row (12) ---
$ws3.Range("A12").Formula = "=SUBTOTAL(9,A7:A11)"
$ws3.Range("B12").Value = "Naman Total"

# --- Saleem rows (13-17) ---
$ws3.Range("A13").Value = 1000
$ws3.Range("B13").Value = "Saleem"
$ws3.Range("A14").Value = 1500
$ws3.Range("B14").Value = "Saleem"
$ws3.Range("A15").Value = 3000
$ws3.Range("B15").Value = "Saleem"
$ws3.Range("A16").Value = 500
$ws3.Range("B16").Value = "Saleem"
$ws3.Range("A17").Value = 1000
$ws3.Range("B17").Value = "Saleem"

$ws3.Range("D13").Value = "Ayaan Min"
$ws3.Range("E13").Formula = "=MINIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Saleem"")"
$ws3.Range("D14").Value = "Ayaan Max"
$ws3.Range("E14").Formula = "=MAXIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Saleem"")"
$ws3.Range("D15").Value = "Ayaan Sum"
$ws3.Range("E15").Formula = "=SUMIFS(`$A`$2:`$A`$18,`$B`$2:`$B`$18,""Saleem"")"

# --- Saleem subtotal row (18) ---
$ws3.Range("A18").Formula = "=SUBTOTAL(9,A13:A17)"
$ws3.Range("B18").Value = "Saleem Total"

# --- Grand total row (19) ---
$ws3.Range("A19").Formula = "=SUBTOTAL(9,A2:A17)"
$ws3.Range("B19").Value = "Grand Total"
